# Gates Demo Final
$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: mark a set of rows' F column TRUE (boolean) ---
$surveyRows = @(16,20,24,28,32,36,52,72)
foreach ($r in $surveyRows) {
    $survey.Cells.Item($r, 6).Value = $true
}

# New selection / scroll position on survey sheet
$survey.Range("F72").Select()
$survey.Application.ActiveWindow.ScrollRow = 47

# Give the survey sheet an explicit (portrait) page setup - this is what
# Excel normally stamps onto pageSetup when a print related property is
# touched via the object model.
$survey.PageSetup.Orientation = 1   # xlPortrait

# --- choices sheet: replace numeric "1" flags with the text "yes" ---
$choicesRows = @(2,3,4,5,6)
foreach ($r in $choicesRows) {
    $choices.Cells.Item($r, 2).Value = "yes"
    $choices.Cells.Item($r, 3).Value = "yes"
}

# Make "choices" the active/selected sheet & tab, with C7 selected
$choices.Activate()
$choices.Range("C7").Select()

$wb.Save()
